# Apply the SPECIES_CODES_AGGREGATES edit:
#  - Remove all rows referencing the AG21 aggregate code
#  - Remove the SMA / SKM row
#  - Add a new SMA / SHM row in its place
#  - Leave the final table as A1:B24, with the last row fully selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-to-top so earlier row numbers stay valid.
$rowsToDelete = @(29, 24, 20, 17, 16, 10, 9)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).EntireRow.Delete()
}

# After the deletions above, row 13 is "SMA"/"AG38" (the last remaining SMA
# row) and row 14 is "FAL"/"SKH". Insert a fresh row 14 for SMA/SHM.
$ws.Rows(14).EntireRow.Insert()
$ws.Cells.Item(14, 1).Value = "SMA"
$ws.Cells.Item(14, 2).Value = "SHM"

# Match the saved selection state from the edit: the whole last row selected.
$ws.Rows("24:24").EntireRow.Select()
